$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'92.167.17"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.095.46"
$ws.Range("E3").Value = "  -1.89%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'233.03"
$ws.Range("E5").Value = "  -2.84%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'610.70"
$ws.Range("E6").Value = "  -1.55%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -3.64%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -0.78%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.04%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "'3.090.11"
$ws.Range("E10").Value = "  -2.10%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.773"
$ws.Range("E11").Value = "  +4.20%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -3.97%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -4.34%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "'91.927.47"
$ws.Range("E14").Value = "  +0.71%  "

# Row 15 - was Avalanche, now Toncoin
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").Value = "'5.37"
$ws.Range("E15").Value = "  -3.93%  "

# Row 16 - was Toncoin, now Avalanche
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'33.42"
$ws.Range("E16").Value = "  -4.67%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "'3.675.80"
$ws.Range("E17").Value = "  -1.52%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.093.02"
$ws.Range("E18").Value = "  -2.26%  "

# Row 19 - SuiNetwork
$ws.Range("D19").Value = "'3.82"
$ws.Range("E19").Value = "  +2.82%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'14.31"
$ws.Range("E20").Value = "  -4.48%  "

# Row 21 - Polkadot
$ws.Range("D21").Value = "'5.71"
$ws.Range("E21").Value = "  -3.42%  "

# Row 22 - was BitcoinCash, now PEPE
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").Value = "'0.0000197"
$ws.Range("E22").Value = "  -2.31%  "

# Row 23 - was PEPE, now BitcoinCash
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'433.11"
$ws.Range("E23").Value = "  -5.30%  "

# Row 24 - Uniswap
$ws.Range("D24").Value = "'9.02"
$ws.Range("E24").Value = "  -1.70%  "

# Row 25 - NEARProtocol
$ws.Range("D25").Value = "'5.54"
$ws.Range("E25").Value = "  -6.11%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "'84.95"
$ws.Range("E26").Value = "  -4.27%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "'11.27"
$ws.Range("E27").Value = "  -4.44%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "'3.259.11"
$ws.Range("E28").Value = "  -1.82%  "

# Row 29 - Dai
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30 - Cronos
$ws.Range("E30").Value = "  +5.50%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "'0.124"
$ws.Range("E31").Value = "  -15.64%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "'0.231"
$ws.Range("E32").Value = "  +0.49%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'9.05"
$ws.Range("E33").Value = "  -3.20%  "

# Row 34 - Binance-PegBSC-USD
$ws.Range("E34").Value = "  -40.19%  "

# Row 35 - RenderToken
$ws.Range("D35").Value = "'7.86"
$ws.Range("E35").Value = "  +5.30%  "

# Row 36 - Kaspa
$ws.Range("D36").Value = "'0.154"
$ws.Range("E36").Value = "  -12.79%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "'25.36"
$ws.Range("E37").Value = "  -3.65%  "

# Row 38 - MantraDAO
$ws.Range("D38").Value = "'3.88"
$ws.Range("E38").Value = "  +0.33%  "

# Row 39 - PancakeSwap
$ws.Range("E39").Value = "  -3.73%  "

# Row 40 - WhiteBITCoin
$ws.Range("D40").Value = "'23.83"
$ws.Range("E40").Value = "  +7.65%  "

# Row 41 - Fetch.AI
$ws.Range("D41").Value = "'1.26"
$ws.Range("E41").Value = "  -4.56%  "

# Row 42 - was Bittensor, now PolygonEcosystemToken
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.432"
$ws.Range("E42").Value = "  -2.85%  "

# Row 43 - was PolygonEcosystemToken, now Bittensor
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'463.58"
$ws.Range("E43").Value = "  -5.64%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  -2.85%  "

# Row 46 - Monero
$ws.Range("D46").Value = "'159.32"
$ws.Range("E46").Value = "  +1.70%  "

# Row 47 - ARBITRUM
$ws.Range("D47").Value = "'0.674"
$ws.Range("E47").Value = "  -5.00%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  -5.40%  "

# Row 49 - was ImmutableX, now OKB
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'43.75"
$ws.Range("E49").Value = "  -0.64%  "

# Row 50 - was OKB, now ImmutableX
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.31"
$ws.Range("E50").Value = "  -3.30%  "

# Row 51 - was VeChain, now FirstDigitalUSD
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'0.997"
$ws.Range("E51").Value = "  -0.08%  "
